# Update Leve profit-tracking cells across all 8 crafting-class sheets
# (currentAveragePrice / NQ / HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ columns H:N)
# with refreshed market-board data from the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 3
$ws.Range("H3").Value = 500653.5
$ws.Range("J3").Value = 500653.5
$ws.Range("L3").Value = 500653.5
$ws.Range("N3").Value = -500881.5

# Row 102
$ws.Range("H102").Value = 500653.5
$ws.Range("J102").Value = 500653.5
$ws.Range("L102").Value = 500653.5
$ws.Range("N102").Value = -507143.5

# Row 137
$ws.Range("H137").Value = 1046.0952
$ws.Range("I137").Value = 936.625
$ws.Range("K137").Value = 2809.875
$ws.Range("M137").Value = -259.875

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2657632.8
$ws.Range("I32").Value = 2502548.2
$ws.Range("K32").Value = 2502548.2
$ws.Range("M32").Value = -2502261.2

# Row 33
$ws.Range("H33").Value = 9999
$ws.Range("I33").Value = 9999
$ws.Range("K33").Value = 9999
$ws.Range("M33").Value = -9670

# Row 45
$ws.Range("H45").Value = 821
$ws.Range("J45").Value = 1014
$ws.Range("L45").Value = 1014
$ws.Range("N45").Value = -1768

# Row 97
$ws.Range("H97").Value = 549.2308
$ws.Range("I97").Value = 563.1
$ws.Range("K97").Value = 563.1
$ws.Range("M97").Value = -67.10000000000002

# Row 122
$ws.Range("H122").Value = 2900
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 1353.1818
$ws.Range("I86").Value = 1388.5
$ws.Range("K86").Value = 1388.5
$ws.Range("M86").Value = -265.5

# Row 89
$ws.Range("H89").Value = 1353.1818
$ws.Range("I89").Value = 1388.5
$ws.Range("K89").Value = 6942.5
$ws.Range("M89").Value = -1326.5

# Row 94
$ws.Range("H94").Value = 480.45456
$ws.Range("I94").Value = 404.22223
$ws.Range("K94").Value = 404.22223
$ws.Range("M94").Value = 46.77776999999998

# Row 105
$ws.Range("H105").Value = 2152
$ws.Range("I105").Value = 1815
$ws.Range("K105").Value = 1815
$ws.Range("M105").Value = -68

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 7999.5
$ws.Range("I16").Value = 5999
$ws.Range("K16").Value = 5999
$ws.Range("M16").Value = -5712

# Row 105
$ws.Range("H105").Value = 2129.818
$ws.Range("I105").Value = 1450.25
$ws.Range("J105").Value = 2945.3
$ws.Range("K105").Value = 1450.25
$ws.Range("L105").Value = 2945.3
$ws.Range("M105").Value = 296.75
$ws.Range("N105").Value = -6439.3

# Row 113
$ws.Range("H113").Value = 7999.5
$ws.Range("I113").Value = 5999
$ws.Range("K113").Value = 5999
$ws.Range("M113").Value = -3829

$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Range("H7").Value = 1000.8333
$ws.Range("J7").Value = 2550
$ws.Range("L7").Value = 7650
$ws.Range("N7").Value = -7874

# Row 50
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("M50").ClearContents()
$ws.Range("N50").ClearContents()

# Row 53
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("M53").ClearContents()
$ws.Range("N53").ClearContents()

# Row 121
$ws.Range("H121").Value = 2457.3333
$ws.Range("I121").Value = 661.6667
$ws.Range("J121").Value = 3355.1667
$ws.Range("K121").Value = 1985.0001
$ws.Range("L121").Value = 10065.5001
$ws.Range("M121").Value = -675.0001
$ws.Range("N121").Value = -12685.5001

# Row 138
$ws.Range("H138").Value = 1817.8334
$ws.Range("I138").Value = 1817.8334
$ws.Range("K138").Value = 5453.5002
$ws.Range("M138").Value = -313.5002000000004

$ws = $wb.Worksheets.Item("GSM")
# Row 9
$ws.Range("H9").Value = 599.6
$ws.Range("I9").Value = 97.5
$ws.Range("K9").Value = 97.5
$ws.Range("M9").Value = 72.5

# Row 80
$ws.Range("H80").Value = 2924.875
$ws.Range("I80").Value = 2449.75
$ws.Range("J80").Value = 3400
$ws.Range("K80").Value = 2449.75
$ws.Range("L80").Value = 3400
$ws.Range("M80").Value = -1451.75
$ws.Range("N80").Value = -5396

# Row 83
$ws.Range("H83").Value = 2924.875
$ws.Range("I83").Value = 2449.75
$ws.Range("J83").Value = 3400
$ws.Range("K83").Value = 12248.75
$ws.Range("L83").Value = 17000
$ws.Range("M83").Value = -7256.75
$ws.Range("N83").Value = -26984

# Row 113
$ws.Range("H113").Value = 924.25
$ws.Range("I113").Value = 849.5
$ws.Range("J113").Value = 999
$ws.Range("K113").Value = 849.5
$ws.Range("L113").Value = 999
$ws.Range("M113").Value = 1320.5
$ws.Range("N113").Value = -5339

# Row 126
$ws.Range("H126").Value = 2000
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -3530
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 18
$ws.Range("H18").Value = 4651.25
$ws.Range("I18").Value = 4651.25
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 4651.25
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -4479.25
$ws.Range("N18").ClearContents()

# Row 20
$ws.Range("H20").Value = 10833.333
$ws.Range("J20").Value = 15000
$ws.Range("L20").Value = 15000
$ws.Range("N20").Value = -15452

# Row 22
$ws.Range("H22").Value = 886.75
$ws.Range("I22").Value = 886.75
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 886.75
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -591.75
$ws.Range("N22").ClearContents()

# Row 27
$ws.Range("H27").Value = 886.75
$ws.Range("I27").Value = 886.75
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 886.75
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -779.75
$ws.Range("N27").ClearContents()

# Row 61
$ws.Range("H61").Value = 1665.6666
$ws.Range("I61").Value = 1499.25
$ws.Range("K61").Value = 1499.25
$ws.Range("M61").Value = -1297.25

# Row 63
$ws.Range("H63").Value = 87077
$ws.Range("I63").Value = 87077
$ws.Range("K63").Value = 87077
$ws.Range("M63").Value = -86328

# Row 66
$ws.Range("H66").Value = 87077
$ws.Range("I66").Value = 87077
$ws.Range("K66").Value = 261231
$ws.Range("M66").Value = -257487

# Row 113
$ws.Range("H113").Value = 1665.6666
$ws.Range("I113").Value = 1499.25
$ws.Range("K113").Value = 1499.25
$ws.Range("M113").Value = 670.75

# Row 122
$ws.Range("H122").Value = 5804.724
$ws.Range("I122").Value = 4790.6875
$ws.Range("K122").Value = 14372.0625
$ws.Range("M122").Value = -11922.0625

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 19166.834
$ws.Range("I62").Value = 20000.666
$ws.Range("K62").Value = 20000.666
$ws.Range("M62").Value = -19376.666

# Row 65
$ws.Range("H65").Value = 19166.834
$ws.Range("I65").Value = 20000.666
$ws.Range("K65").Value = 100003.33
$ws.Range("M65").Value = -96883.33

# Row 68
$ws.Range("H68").Value = 34900
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 34900
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 34900
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -36522

# Row 70
$ws.Range("H70").Value = 49797
$ws.Range("I70").Value = 49595
$ws.Range("K70").Value = 49595
$ws.Range("M70").Value = -49280

# Row 71
$ws.Range("H71").Value = 34900
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 34900
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 104700
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -112812

# Row 73
$ws.Range("H73").Value = 49797
$ws.Range("I73").Value = 49595
$ws.Range("K73").Value = 49595
$ws.Range("M73").Value = -48503

# Row 113
$ws.Range("H113").Value = 723.8333
$ws.Range("I113").Value = 697.5
$ws.Range("K113").Value = 2092.5
$ws.Range("M113").Value = 77.5
